$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates
$ws.Range("D2").Value = "58.086.01"
$ws.Range("D3").Value = "2.972.75"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Value = "2.966.74"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Value = "3.467.00"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "2.970.89"
$ws.Range("D19").Value = "58.005.28"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "422.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0998"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.945"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "0.0₃0699"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "378.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "2.696.80"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.242"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.61"
$ws.Range("D50").Style = "Normal"

# Column E (Volume 1h) updates
$ws.Range("E2").Value = "  -4.10%  "
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  -3.26%  "
$ws.Range("E6").Value = "  +5.55%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +3.08%  "
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("E10").Value = "  -2.67%  "
$ws.Range("E11").Value = "  -5.20%  "
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("E17").Value = "  +7.18%  "
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("E19").Value = "  -4.12%  "
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("E22").Value = "  +3.40%  "
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("E30").Value = "  +5.25%  "
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("E33").Value = "  +6.82%  "
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("E37").Value = "  +5.34%  "
$ws.Range("E38").Value = "  -2.50%  "
$ws.Range("E39").Value = "  +3.24%  "
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  -0.06%  "
